# Regenerate the "K" column (column G) values on the save_data sheet.
# These values are produced by an external calculation (std/mean of
# simulated strikes) and are simply re-written here to match the
# regenerated save data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 4
    5  = 6
    6  = 3
    7  = 2
    8  = 4
    9  = 3
    10 = 1
    11 = 4
    12 = 2
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
